$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "T плавления"
$ws.Range("C1").Value = "T кипения"
$ws.Range("D1").Value = "Уд. теплоемкость в т. с."
$ws.Range("E1").Value = "Уд. теплоёмкость в ж. с."
$ws.Range("F1").Value = "Уд. теплота плавления"
$ws.Range("G1").Value = "Уд. теплота парообр."
$ws.Range("H1").Value = "Плотность в т. с."
$ws.Range("I1").Value = "Плотность в ж. с."

$ws.Range("B1:B2").Select()
